$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.691.42"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.17"
$ws.Range("E3").Value = "  -0.94%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.68%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.36"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4692"
$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3916"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.50"
$ws.Range("E9").Value = "  -4.56%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08001"
$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.005"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.82"
$ws.Range("E12").Value = "  -1.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.877.23"
$ws.Range("E13").Value = "  -0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.990"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.251"
$ws.Range("E15").Value = "  +1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.011"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.59"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06749"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001042"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.23"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.727.19"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.475"
$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.312"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.102.65"
$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.89"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("E28").Value = "  -1.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.151"
$ws.Range("E29").Value = "  +1.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.453"
$ws.Range("E30").Value = "  -2.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.76"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9757"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09527"
$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.631"
$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.328"
$ws.Range("E35").Value = "  -0.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.340"
$ws.Range("E36").Value = "  -7.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06053"
$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02236"
$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.198"
$ws.Range("E39").Value = "  -2.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.292"
$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5990"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1888"
$ws.Range("E43").Value = "  -0.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.31"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.249"
$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5664"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.17"
$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.927"
$ws.Range("E48").Value = "  -1.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06758"
$ws.Range("E49").Value = "  -2.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.03"
$ws.Range("E50").Value = "  -1.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.012"
$ws.Range("E51").Value = "  -11.21%  "
